$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format numeric-looking ID columns (D,E,F) as Text so they are stored
# as strings (vehicleId / Unidad / driverId), matching the source export.
$ws.Range("D2:F4").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = '281474991265672-1752780240944'
$ws.Range("B2").Value = 'No Seat Belt'
$ws.Range("C2").Value = '2025-07-17T13:24:00.944'
$ws.Range("D2").Value = '281474991265672'
$ws.Range("E2").Value = '116'
$ws.Range("F2").Value = '52215735'
$ws.Range("G2").Value = 'ADRIAN CARO'
$ws.Range("H2").Value = 20.67372595
$ws.Range("I2").Value = -103.452576949
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991265672/1752780238444/uVIHaBqhQq-camera-video-segment-driver-1752780240944.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMXZOPPIW%2F20250718%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250718T152641Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHcaCXVzLXdlc3QtMiJHMEUCIQCQW%2BsU5eQwEO75aymYL21PIkaJaxRja72JiSIIBSO7XgIgT6cF%2FA3j95xXTypj353hoeYRQCITjtHHdXZz8WVHqYUq5gMIkP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJKsTSpxwlRWrvfc1yq6Ay9XV7SoJR1I67f4EerNBx32EgXgkQ%2BhYmcnw4MlUckI6X37CNfoO%2BhktUsLTIJpOnj5d21rI5MXj3tp%2B2D5zK4T7O8P1XLHdyy%2B8pDcPtVDL7fRScz7nvD7obmNFtq1KP0o4YQwg7d5D22XW0446Qehn5tmmaXsVjdLNyB3dcrTICzXAS1jxxycvTKnmwDsT5yAb7dfM2fUl3SOxjjkvJPfkyPXsQIpT4JfGvhJb8DdBt55yIrxd8A6TK0OdnJQuNnsitOxcqb1e0Le1szxpXFxaUA0XVACIv6QfvqYOZY9og0QvOaeF%2BG4GkGTYJCcLr%2Bg7ulLdkEtyX4zrpTX09oHZUgA36I6sRVJ6BxZPiVZKChI6dxR1d8y5SdbZJCbVyvwpXSOFo5CsxC%2BA0mJ%2FOiKT8dMhFVwtc8eXMvWMOtKKlUl0siVIpvquOtUPDJSz4EdqxqfKVs6EGKe8R7aPMOi7%2BA2vH%2B168SgwGRNfmnFCGvb01tHtbfoq%2F7dG4wEazGFJc38YhXOiN%2FVQSU1gLPyMwmLhtcwKul%2FQvYbyXi9k9h3urYurdbACHVfBvdexNkE5K43Yw7OeMwwl77pwwY6pQGmueriyrz7cvQTOg6bDqiuexD5%2FETe4Lg8AP50Qe%2BpxICk42zDdjTxc8uViaM8X36Mr4smQNg8twyT9rmsU9sbR1h7s9Ls11zzbrEjCS6e16nTrvuFMGFAWSPH8OCeulaoazxfFiAATM2V8%2B7qEzkNCNY0dFggdzSC15nMBuElSx%2BqvQAgWpbYm8LJQCf8Vvw6m8bVkii7lFngUdlmvm4Z7t51bG0%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2018%20Jul%202025%2023%3A26%3A41%20GMT&X-Amz-Signature=8f6e27cacfbd31982be9415b794721ab6d782fa59c62cc1c6e9ff8d97334e7be'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991265672/1752780238444/YbdfnFXsPB-camera-video-segment-1752780240944.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMXZOPPIW%2F20250718%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250718T152641Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHcaCXVzLXdlc3QtMiJHMEUCIQCQW%2BsU5eQwEO75aymYL21PIkaJaxRja72JiSIIBSO7XgIgT6cF%2FA3j95xXTypj353hoeYRQCITjtHHdXZz8WVHqYUq5gMIkP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJKsTSpxwlRWrvfc1yq6Ay9XV7SoJR1I67f4EerNBx32EgXgkQ%2BhYmcnw4MlUckI6X37CNfoO%2BhktUsLTIJpOnj5d21rI5MXj3tp%2B2D5zK4T7O8P1XLHdyy%2B8pDcPtVDL7fRScz7nvD7obmNFtq1KP0o4YQwg7d5D22XW0446Qehn5tmmaXsVjdLNyB3dcrTICzXAS1jxxycvTKnmwDsT5yAb7dfM2fUl3SOxjjkvJPfkyPXsQIpT4JfGvhJb8DdBt55yIrxd8A6TK0OdnJQuNnsitOxcqb1e0Le1szxpXFxaUA0XVACIv6QfvqYOZY9og0QvOaeF%2BG4GkGTYJCcLr%2Bg7ulLdkEtyX4zrpTX09oHZUgA36I6sRVJ6BxZPiVZKChI6dxR1d8y5SdbZJCbVyvwpXSOFo5CsxC%2BA0mJ%2FOiKT8dMhFVwtc8eXMvWMOtKKlUl0siVIpvquOtUPDJSz4EdqxqfKVs6EGKe8R7aPMOi7%2BA2vH%2B168SgwGRNfmnFCGvb01tHtbfoq%2F7dG4wEazGFJc38YhXOiN%2FVQSU1gLPyMwmLhtcwKul%2FQvYbyXi9k9h3urYurdbACHVfBvdexNkE5K43Yw7OeMwwl77pwwY6pQGmueriyrz7cvQTOg6bDqiuexD5%2FETe4Lg8AP50Qe%2BpxICk42zDdjTxc8uViaM8X36Mr4smQNg8twyT9rmsU9sbR1h7s9Ls11zzbrEjCS6e16nTrvuFMGFAWSPH8OCeulaoazxfFiAATM2V8%2B7qEzkNCNY0dFggdzSC15nMBuElSx%2BqvQAgWpbYm8LJQCf8Vvw6m8bVkii7lFngUdlmvm4Z7t51bG0%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2018%20Jul%202025%2023%3A26%3A41%20GMT&X-Amz-Signature=609302f8cc927faf5073a6b106fcb20c0a1b207480aa5d17db73fb5a1097262c'

# Row 3
$ws.Range("A3").Value = '281474990870452-1752775214763'
$ws.Range("B3").Value = 'Harsh Brake'
$ws.Range("C3").Value = '2025-07-17T12:00:14.763'
$ws.Range("D3").Value = '281474990870452'
$ws.Range("E3").Value = '110'
$ws.Range("F3").Value = '52215670'
$ws.Range("G3").Value = 'ALEJANDRO LARA'
$ws.Range("H3").Value = 20.56563343
$ws.Range("I3").Value = -103.34548504
$ws.Range("J3").Value = 0.5618281364440918
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990870452/1752775209763/sLycc0xtnj-camera-video-segment-driver-1752775214763.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMXZOPPIW%2F20250718%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250718T152641Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHcaCXVzLXdlc3QtMiJHMEUCIQCQW%2BsU5eQwEO75aymYL21PIkaJaxRja72JiSIIBSO7XgIgT6cF%2FA3j95xXTypj353hoeYRQCITjtHHdXZz8WVHqYUq5gMIkP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJKsTSpxwlRWrvfc1yq6Ay9XV7SoJR1I67f4EerNBx32EgXgkQ%2BhYmcnw4MlUckI6X37CNfoO%2BhktUsLTIJpOnj5d21rI5MXj3tp%2B2D5zK4T7O8P1XLHdyy%2B8pDcPtVDL7fRScz7nvD7obmNFtq1KP0o4YQwg7d5D22XW0446Qehn5tmmaXsVjdLNyB3dcrTICzXAS1jxxycvTKnmwDsT5yAb7dfM2fUl3SOxjjkvJPfkyPXsQIpT4JfGvhJb8DdBt55yIrxd8A6TK0OdnJQuNnsitOxcqb1e0Le1szxpXFxaUA0XVACIv6QfvqYOZY9og0QvOaeF%2BG4GkGTYJCcLr%2Bg7ulLdkEtyX4zrpTX09oHZUgA36I6sRVJ6BxZPiVZKChI6dxR1d8y5SdbZJCbVyvwpXSOFo5CsxC%2BA0mJ%2FOiKT8dMhFVwtc8eXMvWMOtKKlUl0siVIpvquOtUPDJSz4EdqxqfKVs6EGKe8R7aPMOi7%2BA2vH%2B168SgwGRNfmnFCGvb01tHtbfoq%2F7dG4wEazGFJc38YhXOiN%2FVQSU1gLPyMwmLhtcwKul%2FQvYbyXi9k9h3urYurdbACHVfBvdexNkE5K43Yw7OeMwwl77pwwY6pQGmueriyrz7cvQTOg6bDqiuexD5%2FETe4Lg8AP50Qe%2BpxICk42zDdjTxc8uViaM8X36Mr4smQNg8twyT9rmsU9sbR1h7s9Ls11zzbrEjCS6e16nTrvuFMGFAWSPH8OCeulaoazxfFiAATM2V8%2B7qEzkNCNY0dFggdzSC15nMBuElSx%2BqvQAgWpbYm8LJQCf8Vvw6m8bVkii7lFngUdlmvm4Z7t51bG0%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2018%20Jul%202025%2023%3A26%3A41%20GMT&X-Amz-Signature=76f406e2e227e6a1439d8ee131b070319d80a34320bc639f3072217f492366f2'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990870452/1752775209763/8y3NajMKj1-camera-video-segment-1752775214763.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMXZOPPIW%2F20250718%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250718T152641Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHcaCXVzLXdlc3QtMiJHMEUCIQCQW%2BsU5eQwEO75aymYL21PIkaJaxRja72JiSIIBSO7XgIgT6cF%2FA3j95xXTypj353hoeYRQCITjtHHdXZz8WVHqYUq5gMIkP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJKsTSpxwlRWrvfc1yq6Ay9XV7SoJR1I67f4EerNBx32EgXgkQ%2BhYmcnw4MlUckI6X37CNfoO%2BhktUsLTIJpOnj5d21rI5MXj3tp%2B2D5zK4T7O8P1XLHdyy%2B8pDcPtVDL7fRScz7nvD7obmNFtq1KP0o4YQwg7d5D22XW0446Qehn5tmmaXsVjdLNyB3dcrTICzXAS1jxxycvTKnmwDsT5yAb7dfM2fUl3SOxjjkvJPfkyPXsQIpT4JfGvhJb8DdBt55yIrxd8A6TK0OdnJQuNnsitOxcqb1e0Le1szxpXFxaUA0XVACIv6QfvqYOZY9og0QvOaeF%2BG4GkGTYJCcLr%2Bg7ulLdkEtyX4zrpTX09oHZUgA36I6sRVJ6BxZPiVZKChI6dxR1d8y5SdbZJCbVyvwpXSOFo5CsxC%2BA0mJ%2FOiKT8dMhFVwtc8eXMvWMOtKKlUl0siVIpvquOtUPDJSz4EdqxqfKVs6EGKe8R7aPMOi7%2BA2vH%2B168SgwGRNfmnFCGvb01tHtbfoq%2F7dG4wEazGFJc38YhXOiN%2FVQSU1gLPyMwmLhtcwKul%2FQvYbyXi9k9h3urYurdbACHVfBvdexNkE5K43Yw7OeMwwl77pwwY6pQGmueriyrz7cvQTOg6bDqiuexD5%2FETe4Lg8AP50Qe%2BpxICk42zDdjTxc8uViaM8X36Mr4smQNg8twyT9rmsU9sbR1h7s9Ls11zzbrEjCS6e16nTrvuFMGFAWSPH8OCeulaoazxfFiAATM2V8%2B7qEzkNCNY0dFggdzSC15nMBuElSx%2BqvQAgWpbYm8LJQCf8Vvw6m8bVkii7lFngUdlmvm4Z7t51bG0%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2018%20Jul%202025%2023%3A26%3A41%20GMT&X-Amz-Signature=0b42f6d625f7ca430238f4600cbde0d4ab8931b24bd0db055b41de5036e882ab'

# Row 4
$ws.Range("A4").Value = '281474991205821-1752773949655'
$ws.Range("B4").Value = 'No Seat Belt'
$ws.Range("C4").Value = '2025-07-17T11:39:09.655'
$ws.Range("D4").Value = '281474991205821'
$ws.Range("E4").Value = '148'
$ws.Range("F4").Value = '51834015'
$ws.Range("G4").Value = 'LUIS IBARRA'
$ws.Range("H4").Value = 20.58452888
$ws.Range("I4").Value = -103.44697295
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1752773947155/MV95RChgaO-camera-video-segment-driver-1752773949655.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMXZOPPIW%2F20250718%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250718T152641Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHcaCXVzLXdlc3QtMiJHMEUCIQCQW%2BsU5eQwEO75aymYL21PIkaJaxRja72JiSIIBSO7XgIgT6cF%2FA3j95xXTypj353hoeYRQCITjtHHdXZz8WVHqYUq5gMIkP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJKsTSpxwlRWrvfc1yq6Ay9XV7SoJR1I67f4EerNBx32EgXgkQ%2BhYmcnw4MlUckI6X37CNfoO%2BhktUsLTIJpOnj5d21rI5MXj3tp%2B2D5zK4T7O8P1XLHdyy%2B8pDcPtVDL7fRScz7nvD7obmNFtq1KP0o4YQwg7d5D22XW0446Qehn5tmmaXsVjdLNyB3dcrTICzXAS1jxxycvTKnmwDsT5yAb7dfM2fUl3SOxjjkvJPfkyPXsQIpT4JfGvhJb8DdBt55yIrxd8A6TK0OdnJQuNnsitOxcqb1e0Le1szxpXFxaUA0XVACIv6QfvqYOZY9og0QvOaeF%2BG4GkGTYJCcLr%2Bg7ulLdkEtyX4zrpTX09oHZUgA36I6sRVJ6BxZPiVZKChI6dxR1d8y5SdbZJCbVyvwpXSOFo5CsxC%2BA0mJ%2FOiKT8dMhFVwtc8eXMvWMOtKKlUl0siVIpvquOtUPDJSz4EdqxqfKVs6EGKe8R7aPMOi7%2BA2vH%2B168SgwGRNfmnFCGvb01tHtbfoq%2F7dG4wEazGFJc38YhXOiN%2FVQSU1gLPyMwmLhtcwKul%2FQvYbyXi9k9h3urYurdbACHVfBvdexNkE5K43Yw7OeMwwl77pwwY6pQGmueriyrz7cvQTOg6bDqiuexD5%2FETe4Lg8AP50Qe%2BpxICk42zDdjTxc8uViaM8X36Mr4smQNg8twyT9rmsU9sbR1h7s9Ls11zzbrEjCS6e16nTrvuFMGFAWSPH8OCeulaoazxfFiAATM2V8%2B7qEzkNCNY0dFggdzSC15nMBuElSx%2BqvQAgWpbYm8LJQCf8Vvw6m8bVkii7lFngUdlmvm4Z7t51bG0%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2018%20Jul%202025%2023%3A26%3A41%20GMT&X-Amz-Signature=d6a1a90942016aac9c3843bc217d2e33369a957d5ee57ce618a5dbcf313a0c4b'
$ws.Range("L4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1752773947155/0QU8iU0WPE-camera-video-segment-1752773949655.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMXZOPPIW%2F20250718%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250718T152641Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHcaCXVzLXdlc3QtMiJHMEUCIQCQW%2BsU5eQwEO75aymYL21PIkaJaxRja72JiSIIBSO7XgIgT6cF%2FA3j95xXTypj353hoeYRQCITjtHHdXZz8WVHqYUq5gMIkP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJKsTSpxwlRWrvfc1yq6Ay9XV7SoJR1I67f4EerNBx32EgXgkQ%2BhYmcnw4MlUckI6X37CNfoO%2BhktUsLTIJpOnj5d21rI5MXj3tp%2B2D5zK4T7O8P1XLHdyy%2B8pDcPtVDL7fRScz7nvD7obmNFtq1KP0o4YQwg7d5D22XW0446Qehn5tmmaXsVjdLNyB3dcrTICzXAS1jxxycvTKnmwDsT5yAb7dfM2fUl3SOxjjkvJPfkyPXsQIpT4JfGvhJb8DdBt55yIrxd8A6TK0OdnJQuNnsitOxcqb1e0Le1szxpXFxaUA0XVACIv6QfvqYOZY9og0QvOaeF%2BG4GkGTYJCcLr%2Bg7ulLdkEtyX4zrpTX09oHZUgA36I6sRVJ6BxZPiVZKChI6dxR1d8y5SdbZJCbVyvwpXSOFo5CsxC%2BA0mJ%2FOiKT8dMhFVwtc8eXMvWMOtKKlUl0siVIpvquOtUPDJSz4EdqxqfKVs6EGKe8R7aPMOi7%2BA2vH%2B168SgwGRNfmnFCGvb01tHtbfoq%2F7dG4wEazGFJc38YhXOiN%2FVQSU1gLPyMwmLhtcwKul%2FQvYbyXi9k9h3urYurdbACHVfBvdexNkE5K43Yw7OeMwwl77pwwY6pQGmueriyrz7cvQTOg6bDqiuexD5%2FETe4Lg8AP50Qe%2BpxICk42zDdjTxc8uViaM8X36Mr4smQNg8twyT9rmsU9sbR1h7s9Ls11zzbrEjCS6e16nTrvuFMGFAWSPH8OCeulaoazxfFiAATM2V8%2B7qEzkNCNY0dFggdzSC15nMBuElSx%2BqvQAgWpbYm8LJQCf8Vvw6m8bVkii7lFngUdlmvm4Z7t51bG0%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2018%20Jul%202025%2023%3A26%3A41%20GMT&X-Amz-Signature=87d134fad7781b3df4f9f2ebc3cb6c1cbf844c86fd68b48b2a4fdfdbadc04079'

Write-Host "Updated rows 2-4 with latest event data"